# CHARMS - Updated NV scenario code to reflect latest changes
# Insert a new question row ("Do you currently live in the United States?")
# right before the existing "In which country do you currently live?" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 20; this shifts all following rows down by one.
$ws.Rows.Item(20).Insert()

# Populate the new row's cells.
$ws.Range("A20").Value2 = "Do you currently live in the United States?"
$ws.Range("B20").Value2 = "Yes"

# Apply formatting to match a bold black question / right-aligned black answer.
$ws.Range("A20").Font.Bold = $true
$ws.Range("A20").Font.Color = 0

$ws.Range("B20").Font.Size = 11
$ws.Range("B20").Font.Color = 0
$ws.Range("B20").HorizontalAlignment = -4152

# Row insertion does not automatically re-anchor existing hyperlinks to their
# new shifted cells, so rebuild them pointing at the correct (shifted) cells.
$ws.Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("B28"), "mailto:consent_participant@yopmail.com") | Out-Null
$ws.Range("B28").HorizontalAlignment = -4152
$ws.Hyperlinks.Add($ws.Range("B29"), "mailto:consent_participant@yopmail.com") | Out-Null
$ws.Range("B29").HorizontalAlignment = -4152

# Update the active selection to reflect where the user was working.
$ws.Range("A21").Select() | Out-Null
